# "Update Naming Plots Regression"
#
# The regression-results block (columns AL:BH) on the "PC2_Pu_GHS" sheet
# previously overflowed past the data table (rows 10-18) into extra rows
# 19-23 that only carried the recycled AL:BH regression output and no
# A:AJ data. This edit removes that stray overflow:
#   - clears the AL:BH regression-output cells from rows 10-18 (the
#     underlying A:AJ data rows stay intact)
#   - deletes rows 19-23 entirely (they held only AL:BH leftovers)
# which also shrinks the sheet's used range down to A1:BH18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PC2_Pu_GHS")

# Remove the spilled-over AL:BH regression values from rows 10-18.
$ws.Range("AL10:BH18").ClearContents()

# Remove the extra rows (19-23) that held only leftover AL:BH values.
$ws.Range("A19:BH23").EntireRow.Delete()
